$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data, sorted descending by value, with "Bengali" and "Uzbek" rows removed.
$data = @(
    @("English", 25.71090207884834),
    @("Chinese", 10.93778660398202),
    @("Spanish", 7.471349227563524),
    @("Japanese", 6.043117907646216),
    @("Arabic", 5.107940598555035),
    @("German", 4.99763041772989),
    @("Russian", 3.531511206663652),
    @("Portuguese", 3.514382822998729),
    @("French", 3.178460763836777),
    @("Italian", 2.851601171445842),
    @("Malay-Indonesian", 2.663440280426213),
    @("Korean", 1.653190034331728),
    @("Persian", 1.549393958128511),
    @("Dutch", 1.477618987507371),
    @("Turkish", 1.457296383855115),
    @("Thai", 0.9808520247197557),
    @("Urdu", 0.8581609537058308),
    @("Polish", 0.8533520891228875),
    @("Swedish", 0.4870659442624978),
    @("Vietnamese", 0.4612988467068549)
)

# Write the reordered rows into A2:B21
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Remove the now-unused trailing rows (previously rows 22 and 23)
$ws.Rows.Item(22).Delete() | Out-Null
$ws.Rows.Item(22).Delete() | Out-Null
